# The diff shows a new data row inserted immediately before the existing
# row 138 ("Pepino ensalada" price record dated 2022-01-01 / serial 44603),
# which pushes the previous rows 138-223 down to 139-224 and grows the
# sheet's used range from A1:R223 to A1:R224.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 138; everything below (138-223) shifts
# down to 139-224, carrying its formatting (incl. the date style on column D).
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new price record.
$ws.Range("A138").Value = 4
$ws.Range("B138").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C138").Value = "Los Lagos"
$ws.Range("D138").Value = 44603
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = 100112043
$ws.Range("G138").Value = "Pepino ensalada"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 400
$ws.Range("K138").Value = 18000
$ws.Range("L138").Value = 18000
$ws.Range("M138").Value = 18000
$ws.Range("N138").Value = "$/caja 60 unidades"
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("P138").Value = 300
$ws.Range("Q138").Value = 60
$ws.Range("R138").Value = "Hortaliza"
